# Auto-generated Excel COM-interop script
# Refresh market-derived columns (H:N) across the 8 crafting-job
# profit sheets from the scheduled Universalis price-fetch run.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 100001540
$ws.Range("I100").Value = 125001050
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 125001050
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -125000509
$ws.Range("N100").Value = -4582
# Row 116
$ws.Range("H116").Value = 5924
$ws.Range("I116").Value = 3125
$ws.Range("K116").Value = 3125
$ws.Range("M116").Value = 317
# Row 129
$ws.Range("H129").Value = 1272.6207
$ws.Range("J129").Value = 1509.826
$ws.Range("L129").Value = 4529.478
$ws.Range("N129").Value = -14529.478
# Row 132
$ws.Range("H132").Value = 6199.067
$ws.Range("I132").Value = 6581.6665
$ws.Range("K132").Value = 19744.9995
$ws.Range("M132").Value = -17214.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4071.65
$ws.Range("I61").Value = 3917.5264
$ws.Range("K61").Value = 3917.5264
$ws.Range("M61").Value = -3705.5264
# Row 97
$ws.Range("H97").Value = 166668600
$ws.Range("I97").Value = 3620
$ws.Range("J97").Value = 333333600
$ws.Range("K97").Value = 3620
$ws.Range("L97").Value = 333333600
$ws.Range("M97").Value = -3124
$ws.Range("N97").Value = -333334592
# Row 122
$ws.Range("H122").Value = 3111.9546
$ws.Range("I122").Value = 2466.5625
$ws.Range("J122").Value = 4833
$ws.Range("K122").Value = 7399.6875
$ws.Range("L122").Value = 14499
$ws.Range("M122").Value = -4949.6875
$ws.Range("N122").Value = -19399
# Row 132
$ws.Range("H132").Value = 31132.723
$ws.Range("I132").Value = 3378.5454
$ws.Range("K132").Value = 10135.6362
$ws.Range("M132").Value = -7605.636200000001
# Row 136
$ws.Range("H136").Value = 4071.65
$ws.Range("I136").Value = 3917.5264
$ws.Range("K136").Value = 11752.5792
$ws.Range("M136").Value = -9202.5792

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
# Row 107
$ws.Range("H107").Value = 860.4167
$ws.Range("I107").Value = 905.1111
$ws.Range("K107").Value = 905.1111
$ws.Range("M107").Value = 1014.8889
# Row 134
$ws.Range("H134").Value = 7823.25
$ws.Range("I134").Value = 8352.637000000001
$ws.Range("K134").Value = 25057.911
$ws.Range("M134").Value = -22522.911

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 58
$ws.Range("H58").Value = 26966.8
$ws.Range("I58").Value = 2149.25
$ws.Range("J58").Value = 43511.832
$ws.Range("K58").Value = 2149.25
$ws.Range("L58").Value = 43511.832
$ws.Range("M58").Value = -1946.25
$ws.Range("N58").Value = -43917.832
# Row 62
$ws.Range("H62").Value = 83337160
$ws.Range("J62").Value = 3300
$ws.Range("L62").Value = 3300
$ws.Range("N62").Value = -4548
# Row 65
$ws.Range("H65").Value = 83337160
$ws.Range("J65").Value = 3300
$ws.Range("L65").Value = 16500
$ws.Range("N65").Value = -22740
# Row 68
$ws.Range("H68").Value = 50725
$ws.Range("J68").Value = 50725
$ws.Range("L68").Value = 50725
$ws.Range("N68").Value = -52223
# Row 71
$ws.Range("H71").Value = 50725
$ws.Range("J71").Value = 50725
$ws.Range("L71").Value = 152175
$ws.Range("N71").Value = -159663
# Row 94
$ws.Range("H94").Value = 2839.6924
$ws.Range("I94").Value = 637.2
$ws.Range("J94").Value = 4216.25
$ws.Range("K94").Value = 637.2
$ws.Range("L94").Value = 4216.25
$ws.Range("M94").Value = -186.2
$ws.Range("N94").Value = -5118.25
# Row 107
$ws.Range("H107").Value = 1724.421
$ws.Range("I107").Value = 1598.3334
$ws.Range("J107").Value = 1782.6154
$ws.Range("K107").Value = 1598.3334
$ws.Range("L107").Value = 1782.6154
$ws.Range("M107").Value = 321.6666
$ws.Range("N107").Value = -5622.6154
# Row 132
$ws.Range("H132").Value = 4729.6
$ws.Range("I132").Value = 2322.3333
$ws.Range("K132").Value = 6966.999899999999
$ws.Range("M132").Value = -4436.999899999999
# Row 134
$ws.Range("H134").Value = 1244.6364
$ws.Range("I134").Value = 961.375
$ws.Range("K134").Value = 2884.125
$ws.Range("M134").Value = -349.125
# Row 136
$ws.Range("H136").Value = 26966.8
$ws.Range("I136").Value = 2149.25
$ws.Range("J136").Value = 43511.832
$ws.Range("K136").Value = 6447.75
$ws.Range("L136").Value = 130535.496
$ws.Range("M136").Value = -3897.75
$ws.Range("N136").Value = -135635.496

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 596.6667
$ws.Range("I17").Value = 450
$ws.Range("J17").Value = 670
$ws.Range("K17").Value = 1350
$ws.Range("L17").Value = 2010
$ws.Range("M17").Value = -1181
$ws.Range("N17").Value = -2348
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
# Row 121
$ws.Range("H121").Value = 948.4838999999999
$ws.Range("I121").Value = 501.625
$ws.Range("J121").Value = 1103.9131
$ws.Range("K121").Value = 1504.875
$ws.Range("L121").Value = 3311.7393
$ws.Range("M121").Value = -194.875
$ws.Range("N121").Value = -5931.7393
# Row 131
$ws.Range("H131").Value = 702.44446
$ws.Range("J131").Value = 720.34784
$ws.Range("L131").Value = 2161.04352
$ws.Range("N131").Value = -12241.04352

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 222224640
$ws.Range("I122").Value = 66669564
$ws.Range("J122").Value = 1000000000
$ws.Range("K122").Value = 200008692
$ws.Range("L122").Value = 3000000000
$ws.Range("M122").Value = -200006242
$ws.Range("N122").Value = -3000004900

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 753.3333
$ws.Range("I16").Value = 697.5
$ws.Range("K16").Value = 697.5
$ws.Range("M16").Value = -527.5
# Row 61
$ws.Range("H61").Value = 5491.5835
$ws.Range("I61").Value = 2320
$ws.Range("K61").Value = 2320
$ws.Range("M61").Value = -2118
# Row 68
$ws.Range("H68").Value = 2547.8333
$ws.Range("I68").Value = 1100
$ws.Range("J68").Value = 2837.4
$ws.Range("K68").Value = 1100
$ws.Range("L68").Value = 2837.4
$ws.Range("M68").Value = -351
$ws.Range("N68").Value = -4335.4
# Row 71
$ws.Range("H71").Value = 2547.8333
$ws.Range("I71").Value = 1100
$ws.Range("J71").Value = 2837.4
$ws.Range("K71").Value = 5500
$ws.Range("L71").Value = 14187
$ws.Range("M71").Value = -1756
$ws.Range("N71").Value = -21675
# Row 113
$ws.Range("H113").Value = 5491.5835
$ws.Range("I113").Value = 2320
$ws.Range("K113").Value = 2320
$ws.Range("M113").Value = -150
# Row 132
$ws.Range("H132").Value = 3741.7144
$ws.Range("I132").Value = 2838.8
$ws.Range("K132").Value = 8516.400000000001
$ws.Range("M132").Value = -5986.400000000001
# Row 136
$ws.Range("H136").Value = 2210.7058
$ws.Range("I136").Value = 2143.818
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 6431.454000000001
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = -3881.454000000001
$ws.Range("N136").Value = -12099.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 2980.6667
$ws.Range("I113").Value = 3520.8
$ws.Range("K113").Value = 10562.4
$ws.Range("M113").Value = -8392.400000000001
# Row 132
$ws.Range("H132").Value = 1916.3572
$ws.Range("I132").Value = 1264.8572
$ws.Range("J132").Value = 2567.8572
$ws.Range("K132").Value = 3794.5716
$ws.Range("L132").Value = 7703.571599999999
$ws.Range("M132").Value = -1264.5716
$ws.Range("N132").Value = -12763.5716
# Row 136
$ws.Range("H136").Value = 22224676
$ws.Range("I136").Value = 30304176
$ws.Range("J136").Value = 6050.25
$ws.Range("K136").Value = 90912528
$ws.Range("L136").Value = 18150.75
$ws.Range("N136").Value = -23250.75
